# Auto-generated script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.700.21"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.919.89"
$ws.Range("E3").Value = "  -2.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "374.22"
$ws.Range("E5").Value = "  -2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.61"
$ws.Range("E6").Value = "  -2.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.534"
$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  -3.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").Value = "  -3.24%  "

$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.374.32"
$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("B14").Value = "Uniswap"
$ws.Range("C14").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.11"
$ws.Range("E14").Value = "  +64.30%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "17.94"
$ws.Range("E15").Value = "  -1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.55"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.909.71"
$ws.Range("E17").Value = "  -2.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.989"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "50.725.23"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.98"
$ws.Range("E20").Value = "  -7.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.18"
$ws.Range("E21").Value = "  -4.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.18"
$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.89"
$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  +8.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.79"
$ws.Range("E26").Value = "  -4.94%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").Value = "  -6.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.33"
$ws.Range("E29").Value = "  -1.99%  "

$ws.Range("E30").Value = "  -4.95%  "

$ws.Range("E31").Value = "  -5.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.86"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.05"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.27"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.02"
$ws.Range("E35").Value = "  -4.15%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0428"
$ws.Range("E36").Value = "  -4.32%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.05"
$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.23"
$ws.Range("E40").Value = "  -4.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  -1.27%  "

$ws.Range("E42").Value = "  -6.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.18"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.83"
$ws.Range("E44").Value = "  -3.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("E45").Value = "  +2.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("E46").Value = "  -2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.984.04"
$ws.Range("E48").Value = "  -2.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.257"
$ws.Range("E49").Value = "  -7.03%  "

$ws.Range("E50").Value = "  -6.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.25"
$ws.Range("E51").Value = "  +2.73%  "

Write-Host "cryptos list updated"
